$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 800.1539
$ws.Range("I41").Value = 920
$ws.Range("J41").Value = 725.25
$ws.Range("K41").Value = 920
$ws.Range("L41").Value = 725.25
$ws.Range("M41").Value = -480
$ws.Range("N41").Value = -1605.25
$ws.Range("H125").Value = 2776
$ws.Range("I125").Value = 997.7143
$ws.Range("J125").Value = 9000
$ws.Range("K125").Value = 8979.4287
$ws.Range("L125").Value = 81000
$ws.Range("M125").Value = -6519.4287
$ws.Range("N125").Value = -85920
$ws.Range("H132").Value = 1274.6666
$ws.Range("I132").Value = 1236
$ws.Range("J132").Value = 1410
$ws.Range("K132").Value = 3708
$ws.Range("L132").Value = 4230
$ws.Range("M132").Value = -1178
$ws.Range("N132").Value = -9290
$ws.Range("H137").Value = 1332.9697
$ws.Range("I137").Value = 1147.44
$ws.Range("K137").Value = 3442.32
$ws.Range("M137").Value = -892.3200000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4700.7334
$ws.Range("I61").Value = 4941.893
$ws.Range("J61").Value = 1324.5
$ws.Range("K61").Value = 4941.893
$ws.Range("L61").Value = 1324.5
$ws.Range("M61").Value = -4729.893
$ws.Range("N61").Value = -1748.5
$ws.Range("H74").Value = 4963.6294
$ws.Range("I74").Value = 1229.3334
$ws.Range("K74").Value = 1229.3334
$ws.Range("M74").Value = -355.3334
$ws.Range("H77").Value = 4963.6294
$ws.Range("I77").Value = 1229.3334
$ws.Range("K77").Value = 6146.666999999999
$ws.Range("M77").Value = -1778.666999999999
$ws.Range("H88").Value = 2714.2727
$ws.Range("I88").Value = 2650.5
$ws.Range("J88").Value = 2790.8
$ws.Range("K88").Value = 2650.5
$ws.Range("L88").Value = 2790.8
$ws.Range("M88").Value = -2244.5
$ws.Range("N88").Value = -3602.8
$ws.Range("H91").Value = 2714.2727
$ws.Range("I91").Value = 2650.5
$ws.Range("J91").Value = 2790.8
$ws.Range("K91").Value = 2650.5
$ws.Range("L91").Value = 2790.8
$ws.Range("M91").Value = -1246.5
$ws.Range("N91").Value = -5598.8
$ws.Range("H122").Value = 1351163.2
$ws.Range("I122").Value = 1832599.1
$ws.Range("J122").Value = 3142.8
$ws.Range("K122").Value = 5497797.300000001
$ws.Range("L122").Value = 9428.400000000001
$ws.Range("M122").Value = -5495347.300000001
$ws.Range("N122").Value = -14328.4
$ws.Range("H132").Value = 3995.1292
$ws.Range("I132").Value = 2412.5264
$ws.Range("J132").Value = 6500.9165
$ws.Range("K132").Value = 7237.5792
$ws.Range("L132").Value = 19502.7495
$ws.Range("M132").Value = -4707.5792
$ws.Range("N132").Value = -24562.7495
$ws.Range("H136").Value = 4700.7334
$ws.Range("I136").Value = 4941.893
$ws.Range("J136").Value = 1324.5
$ws.Range("K136").Value = 14825.679
$ws.Range("L136").Value = 3973.5
$ws.Range("M136").Value = -12275.679
$ws.Range("N136").Value = -9073.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 50780
$ws.Range("J59").Value = 50780
$ws.Range("L59").Value = 50780
$ws.Range("N59").Value = -52474
$ws.Range("H86").Value = 10102716
$ws.Range("I86").Value = 17545506
$ws.Range("K86").Value = 17545506
$ws.Range("M86").Value = -17544383
$ws.Range("H89").Value = 10102716
$ws.Range("I89").Value = 17545506
$ws.Range("K89").Value = 87727530
$ws.Range("M89").Value = -87721914
$ws.Range("H105").Value = 2328.318
$ws.Range("I105").Value = 1663.4286
$ws.Range("J105").Value = 3491.875
$ws.Range("K105").Value = 1663.4286
$ws.Range("L105").Value = 3491.875
$ws.Range("M105").Value = 83.57140000000004
$ws.Range("N105").Value = -6985.875
$ws.Range("H134").Value = 6318.846
$ws.Range("I134").Value = 8544.0625
$ws.Range("J134").Value = 2758.5
$ws.Range("K134").Value = 25632.1875
$ws.Range("L134").Value = 8275.5
$ws.Range("M134").Value = -23097.1875
$ws.Range("N134").Value = -13345.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1252.9412
$ws.Range("I16").Value = 1127.2727
$ws.Range("J16").Value = 1483.3334
$ws.Range("K16").Value = 1127.2727
$ws.Range("L16").Value = 1483.3334
$ws.Range("M16").Value = -840.2727
$ws.Range("N16").Value = -2057.3334
$ws.Range("H31").Value = 5906.756
$ws.Range("I31").Value = 2175.5264
$ws.Range("J31").Value = 9129.182000000001
$ws.Range("K31").Value = 2175.5264
$ws.Range("L31").Value = 9129.182000000001
$ws.Range("M31").Value = -1880.5264
$ws.Range("N31").Value = -9719.182000000001
$ws.Range("H34").Value = 5906.756
$ws.Range("I34").Value = 2175.5264
$ws.Range("J34").Value = 9129.182000000001
$ws.Range("K34").Value = 2175.5264
$ws.Range("L34").Value = 9129.182000000001
$ws.Range("M34").Value = -1973.5264
$ws.Range("N34").Value = -9533.182000000001
$ws.Range("H62").Value = 4700.2856
$ws.Range("I62").Value = 5153.846
$ws.Range("J62").Value = 3963.25
$ws.Range("K62").Value = 5153.846
$ws.Range("L62").Value = 3963.25
$ws.Range("M62").Value = -4529.846
$ws.Range("N62").Value = -5211.25
$ws.Range("H65").Value = 4700.2856
$ws.Range("I65").Value = 5153.846
$ws.Range("J65").Value = 3963.25
$ws.Range("K65").Value = 25769.23
$ws.Range("L65").Value = 19816.25
$ws.Range("M65").Value = -22649.23
$ws.Range("N65").Value = -26056.25
$ws.Range("H99").Value = 3257.3845
$ws.Range("I99").Value = 2288.7368
$ws.Range("J99").Value = 5886.5713
$ws.Range("K99").Value = 2288.7368
$ws.Range("L99").Value = 5886.5713
$ws.Range("M99").Value = -790.7368000000001
$ws.Range("N99").Value = -8882.5713
$ws.Range("H105").Value = 1542.2
$ws.Range("I105").Value = 1900
$ws.Range("J105").Value = 1303.6666
$ws.Range("K105").Value = 1900
$ws.Range("L105").Value = 1303.6666
$ws.Range("M105").Value = -153
$ws.Range("N105").Value = -4797.6666
$ws.Range("H113").Value = 1252.9412
$ws.Range("I113").Value = 1127.2727
$ws.Range("J113").Value = 1483.3334
$ws.Range("K113").Value = 1127.2727
$ws.Range("L113").Value = 1483.3334
$ws.Range("M113").Value = 1042.7273
$ws.Range("N113").Value = -5823.3334
$ws.Range("H122").Value = 841.2273
$ws.Range("I122").Value = 628.5714
$ws.Range("J122").Value = 1213.375
$ws.Range("K122").Value = 1885.7142
$ws.Range("L122").Value = 3640.125
$ws.Range("M122").Value = 564.2857999999999
$ws.Range("N122").Value = -8540.125
$ws.Range("H126").Value = 3257.3845
$ws.Range("I126").Value = 2288.7368
$ws.Range("J126").Value = 5886.5713
$ws.Range("K126").Value = 6866.2104
$ws.Range("L126").Value = 17659.7139
$ws.Range("M126").Value = -4396.2104
$ws.Range("N126").Value = -22599.7139
$ws.Range("H132").Value = 2107.5217
$ws.Range("I132").Value = 1949.6097
$ws.Range("J132").Value = 3402.4
$ws.Range("K132").Value = 5848.8291
$ws.Range("L132").Value = 10207.2
$ws.Range("M132").Value = -3318.8291
$ws.Range("N132").Value = -15267.2
$ws.Range("H134").Value = 5472
$ws.Range("I134").Value = 6949
$ws.Range("J134").Value = 1533.3334
$ws.Range("K134").Value = 20847
$ws.Range("L134").Value = 4600.0002
$ws.Range("M134").Value = -18312
$ws.Range("N134").Value = -9670.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 274008.6
$ws.Range("J5").Value = 669045
$ws.Range("L5").Value = 2007135
$ws.Range("N5").Value = -2007359
$ws.Range("H107").Value = 768.5714
$ws.Range("I107").Value = 780
$ws.Range("J107").Value = 766.6667
$ws.Range("K107").Value = 2340
$ws.Range("L107").Value = 2300.0001
$ws.Range("M107").Value = -420
$ws.Range("N107").Value = -6140.0001
$ws.Range("H135").Value = 274008.6
$ws.Range("J135").Value = 669045
$ws.Range("L135").Value = 6021405
$ws.Range("N135").Value = -6026475

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 58824750
$ws.Range("I113").Value = 166667400
$ws.Range("K113").Value = 166667400
$ws.Range("M113").Value = -166665230

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 57633.5
$ws.Range("I7").Value = 112833.664
$ws.Range("J7").Value = 2433.3333
$ws.Range("K7").Value = 112833.664
$ws.Range("L7").Value = 2433.3333
$ws.Range("M7").Value = -112721.664
$ws.Range("N7").Value = -2657.3333
$ws.Range("H61").Value = 3368.75
$ws.Range("I61").Value = 2750
$ws.Range("J61").Value = 3987.5
$ws.Range("K61").Value = 2750
$ws.Range("L61").Value = 3987.5
$ws.Range("M61").Value = -2548
$ws.Range("N61").Value = -4391.5
$ws.Range("H69").Value = 100000
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 100000
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H113").Value = 3368.75
$ws.Range("I113").Value = 2750
$ws.Range("J113").Value = 3987.5
$ws.Range("K113").Value = 2750
$ws.Range("L113").Value = 3987.5
$ws.Range("M113").Value = -580
$ws.Range("N113").Value = -8327.5
$ws.Range("H126").Value = 57633.5
$ws.Range("I126").Value = 112833.664
$ws.Range("J126").Value = 2433.3333
$ws.Range("K126").Value = 338500.992
$ws.Range("L126").Value = 7299.999899999999
$ws.Range("M126").Value = -336030.992
$ws.Range("N126").Value = -12239.9999
$ws.Range("H136").Value = 9852.862999999999
$ws.Range("I136").Value = 10546.385
$ws.Range("J136").Value = 8851.111000000001
$ws.Range("K136").Value = 31639.155
$ws.Range("L136").Value = 26553.333
$ws.Range("M136").Value = -29089.155
$ws.Range("N136").Value = -31653.333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 682.73914
$ws.Range("I113").Value = 535.6
$ws.Range("J113").Value = 958.625
$ws.Range("K113").Value = 1606.8
$ws.Range("L113").Value = 2875.875
$ws.Range("M113").Value = 563.1999999999998
$ws.Range("N113").Value = -7215.875
$ws.Range("H122").Value = 1089.4546
$ws.Range("I122").Value = 1088.4
$ws.Range("J122").Value = 1100
$ws.Range("K122").Value = 3265.2
$ws.Range("L122").Value = 3300
$ws.Range("M122").Value = -815.2000000000003
$ws.Range("N122").Value = -8200
$ws.Range("H132").Value = 1741
$ws.Range("I132").Value = 1210.9412
$ws.Range("J132").Value = 2491.9167
$ws.Range("K132").Value = 3632.8236
$ws.Range("L132").Value = 7475.750100000001
$ws.Range("M132").Value = -1102.8236
$ws.Range("N132").Value = -12535.7501
